$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins / Losses / Ties, matching the header style used
# by the rest of row 1 (copy style from AC1 which is the last existing header cell).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1:AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record values (67 wins, 94 losses, 0 ties) for every data row (2-49)
$lastRow = 49
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 67   # AD
    $ws.Cells.Item($r, 31).Value = 94   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
